$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acessos")
$ws2 = $wb.Worksheets.Item("Categorias")

# Turn off the existing autofilter so we can redefine its range later
$ws.AutoFilterMode = $false

$ws.Range("C2").Value = "BOX"

$ws.Range("B3").Value = "Comercial"
$ws.Range("C3").Value = "3PL"

$ws.Range("A4").Value = "Teste"

$ws.Range("A5").Value = "Teste1"
$ws.Range("B5").Value = "Técnica"
$ws.Range("C5").Value = "FLYERS"

$ws.Range("A6").Value = "Teste2"
$ws.Range("C6").Value = "CALL CENTER"
$ws.Range("C6").Font.Color = 0

$ws.Range("D1").Value = "Senha"

$ws.Range("A7").Value = "Teste3"
$ws.Range("B7").Value = "ESG"
$ws.Range("C7").Value = "3PL"

$ws.Range("A8").Value = "brunojeliel"
$ws.Range("B8").Value = "ESG"
$ws.Range("C8").Value = "FLYERS"

# Re-create the autofilter over the expanded range A1:D1
$ws.Range("A1:D1").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Acessos!_FilterDatabase") {
        $n.RefersTo = "=Acessos!`$A`$1:`$D`$1"
    }
}

$ws2.Range("A73").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 46

$ws.Activate()
$ws.Range("C14").Select() | Out-Null
